# Delete row 420 ("「おやすみザイナ」..." post) entirely.
# This shifts all subsequent rows up by one, matching the target diff
# where row 420 is removed and rows 421-621 are renumbered to 420-620.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(420).Delete()
